$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header title text (period after "5.6.1.1" removed)
$ws.Range("B1").Value = "5.6.1.1 Доля замужних женщин и сексуально активных не замужних женщин в возрасте 15-49 лет, которые были осведомлены о соврменном методе контрацепции"

# Update the "urban"/"rural" row labels to their expanded forms
$ws.Range("A6").Value = "Шаар жерлери"
$ws.Range("B6").Value = "Городские поселения"
$ws.Range("C6").Value = "City"

$ws.Range("A7").Value = "Айыл аймагы"
$ws.Range("B7").Value = "Сельская местность"
$ws.Range("C7").Value = "Village"

# Update the active selection shown in the sheet view
$ws.Activate()
$ws.Range("A6:C7").Select()
